# Update the Markov transition-probability matrix on Sheet1 with refreshed
# values reflecting additional simulated games (see commit message:
# "added more games, sped up simulate game logic, and drafted optimization logic").
# Each row (a starting state in column A) is a probability distribution over the
# possible next states (columns B:S) and still sums to 1 after the update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1985559566787004
$ws.Range("C2").Value = 0.5740072202166066
$ws.Range("J2").Value = 0.007220216606498195
$ws.Range("P2").Value = 0.1552346570397112
$ws.Range("S2").Value = 0.06498194945848375
$ws.Range("B3").Value = 0.01257861635220126
$ws.Range("C3").Value = 0.006289308176100629
$ws.Range("J3").Value = 0.03144654088050314
$ws.Range("P3").Value = 0.710691823899371
$ws.Range("S3").Value = 0.2389937106918239
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("B6").Value = 0.04629629629629629
$ws.Range("D6").Value = 0.01388888888888889
$ws.Range("F6").Value = 0.07870370370370371
$ws.Range("J6").Value = 0.2685185185185185
$ws.Range("O6").Value = 0.009259259259259259
$ws.Range("Q6").Value = 0.1712962962962963
$ws.Range("R6").Value = 0.08796296296296297
$ws.Range("S6").Value = 0.3240740740740741
$ws.Range("B7").Value = 0.07291666666666667
$ws.Range("D7").Value = 0.02083333333333333
$ws.Range("E7").Value = 0.005208333333333333
$ws.Range("F7").Value = 0.0625
$ws.Range("J7").Value = 0.1770833333333333
$ws.Range("O7").Value = 0.015625
$ws.Range("Q7").Value = 0.15625
$ws.Range("R7").Value = 0.109375
$ws.Range("S7").Value = 0.3802083333333333
$ws.Range("B8").Value = 0.08646616541353383
$ws.Range("D8").Value = 0.007518796992481203
$ws.Range("F8").Value = 0.07330827067669173
$ws.Range("J8").Value = 0.1221804511278195
$ws.Range("O8").Value = 0.02067669172932331
$ws.Range("Q8").Value = 0.1616541353383459
$ws.Range("R8").Value = 0.09398496240601503
$ws.Range("S8").Value = 0.4342105263157895
$ws.Range("B9").Value = 0.09859154929577464
$ws.Range("D9").Value = 0.01408450704225352
$ws.Range("F9").Value = 0.05164319248826291
$ws.Range("J9").Value = 0.1173708920187793
$ws.Range("O9").Value = 0.02816901408450704
$ws.Range("Q9").Value = 0.1737089201877934
$ws.Range("R9").Value = 0.08450704225352113
$ws.Range("S9").Value = 0.431924882629108
$ws.Range("B10").Value = 0.09738717339667459
$ws.Range("D10").Value = 0.02058590657165479
$ws.Range("F10").Value = 0.05700712589073634
$ws.Range("J10").Value = 0.1258907363420428
$ws.Range("O10").Value = 0.01108471892319873
$ws.Range("Q10").Value = 0.2193190815518606
$ws.Range("R10").Value = 0.09026128266033254
$ws.Range("S10").Value = 0.3784639746634996
$ws.Range("G11").Value = 0.1569230769230769
$ws.Range("J11").Value = 0.08923076923076922
$ws.Range("K11").Value = 0.2307692307692308
$ws.Range("L11").Value = 0.4830769230769231
$ws.Range("S11").Value = 0.04
$ws.Range("G12").Value = 0.703030303030303
$ws.Range("J12").Value = 0.1575757575757576
$ws.Range("K12").Value = 0.02424242424242424
$ws.Range("L12").Value = 0.06060606060606061
$ws.Range("S12").Value = 0.05454545454545454
$ws.Range("G13").Value = 0.6829268292682927
$ws.Range("J13").Value = 0.2439024390243902
$ws.Range("S13").Value = 0.07317073170731707
$ws.Range("F15").Value = 0.01538461538461539
$ws.Range("H15").Value = 0.1641025641025641
$ws.Range("I15").Value = 0.03076923076923077
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.06153846153846154
$ws.Range("M15").Value = 0.01025641025641026
$ws.Range("O15").Value = 0.05641025641025641
$ws.Range("S15").Value = 0.2615384615384616
$ws.Range("F16").Value = 0.01694915254237288
$ws.Range("H16").Value = 0.2259887005649718
$ws.Range("I16").Value = 0.096045197740113
$ws.Range("J16").Value = 0.4180790960451977
$ws.Range("K16").Value = 0.096045197740113
$ws.Range("M16").Value = 0.005649717514124294
$ws.Range("O16").Value = 0.02259887005649718
$ws.Range("S16").Value = 0.1186440677966102
$ws.Range("F17").Value = 0.01939655172413793
$ws.Range("H17").Value = 0.2004310344827586
$ws.Range("I17").Value = 0.08836206896551724
$ws.Range("J17").Value = 0.3685344827586207
$ws.Range("K17").Value = 0.09913793103448276
$ws.Range("M17").Value = 0.01939655172413793
$ws.Range("O17").Value = 0.0625
$ws.Range("S17").Value = 0.1422413793103448
$ws.Range("F18").Value = 0.009009009009009009
$ws.Range("H18").Value = 0.2027027027027027
$ws.Range("I18").Value = 0.0990990990990991
$ws.Range("J18").Value = 0.3828828828828829
$ws.Range("K18").Value = 0.1036036036036036
$ws.Range("M18").Value = 0.004504504504504504
$ws.Range("N18").Value = 0.004504504504504504
$ws.Range("O18").Value = 0.04954954954954955
$ws.Range("S18").Value = 0.1441441441441441
$ws.Range("F19").Value = 0.01388888888888889
$ws.Range("H19").Value = 0.2390350877192982
$ws.Range("I19").Value = 0.0935672514619883
$ws.Range("J19").Value = 0.3340643274853801
$ws.Range("K19").Value = 0.1052631578947368
$ws.Range("M19").Value = 0.02192982456140351
$ws.Range("N19").Value = 0.002192982456140351
$ws.Range("O19").Value = 0.05628654970760234
$ws.Range("S19").Value = 0.1337719298245614
